$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header formatting (bold, bordered, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 94
    $ws.Cells.Item($r, 31).Value = 68
    $ws.Cells.Item($r, 32).Value = 0
}
